$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")
$ws.Activate()

# Insert a new row above row 12 (current "status" row) so that the
# existing status / status_reason / date_created rows shift down by one,
# then populate the freed row with the new "assigned_item_pack_code" field
# (needed for the CTP use case, alongside entitlements/authorizations).
$ws.Rows.Item(12).Insert()

# Copy the formatting from the row above (same font/fill as every other
# data row) onto the freshly inserted row before writing its values.
$ws.Range("A11:B11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "string"
$ws.Range("B12").Value = "assigned_item_pack_code"

$ws.Range("C14").Select()
